$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 599.6087
$ws.Cells.Item(92, 9).Value = 535.94116
$ws.Cells.Item(92, 10).Value = 780
$ws.Cells.Item(92, 11).Value = 535.94116
$ws.Cells.Item(92, 12).Value = 780
$ws.Cells.Item(92, 13).Value = 712.05884
$ws.Cells.Item(92, 14).Value = -3276
$ws.Cells.Item(112, 8).Value = 2584988.5
$ws.Cells.Item(112, 9).Value = 787.5
$ws.Cells.Item(112, 10).Value = 2850034.8
$ws.Cells.Item(112, 11).Value = 2362.5
$ws.Cells.Item(112, 12).Value = 8550104.399999999
$ws.Cells.Item(112, 13).Value = -1254.5
$ws.Cells.Item(112, 14).Value = -8552320.399999999
$ws.Cells.Item(129, 8).Value = 170510.25
$ws.Cells.Item(129, 10).Value = 182884.64
$ws.Cells.Item(129, 12).Value = 548653.92
$ws.Cells.Item(129, 14).Value = -558653.92
$ws.Cells.Item(132, 8).Value = 2643.1
$ws.Cells.Item(132, 9).Value = 2482.9744
$ws.Cells.Item(132, 11).Value = 7448.9232
$ws.Cells.Item(132, 13).Value = -4918.9232
$ws.Cells.Item(138, 8).Value = 33336652
$ws.Cells.Item(138, 10).Value = 3336.652
$ws.Cells.Item(138, 12).Value = 10009.956
$ws.Cells.Item(138, 14).Value = -20289.956

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1865
$ws.Cells.Item(2, 9).Value = 1715.5834
$ws.Cells.Item(2, 11).Value = 1715.5834
$ws.Cells.Item(2, 13).Value = -1602.5834
$ws.Cells.Item(32, 8).Value = 5340.705
$ws.Cells.Item(32, 9).Value = 4541.518
$ws.Cells.Item(32, 11).Value = 4541.518
$ws.Cells.Item(32, 13).Value = -4254.518
$ws.Cells.Item(61, 8).Value = 3770.5881
$ws.Cells.Item(61, 9).Value = 3846.875
$ws.Cells.Item(61, 10).Value = 2550
$ws.Cells.Item(61, 11).Value = 3846.875
$ws.Cells.Item(61, 12).Value = 2550
$ws.Cells.Item(61, 13).Value = -3634.875
$ws.Cells.Item(61, 14).Value = -2974
$ws.Cells.Item(97, 8).Value = 100001230
$ws.Cells.Item(97, 9).Value = 1162.5
$ws.Cells.Item(97, 11).Value = 1162.5
$ws.Cells.Item(97, 13).Value = -666.5
$ws.Cells.Item(102, 8).Value = 2000
$ws.Cells.Item(102, 9).Value = 2000
$ws.Cells.Item(102, 11).Value = 2000
$ws.Cells.Item(102, 13).Value = -378
$ws.Cells.Item(116, 8).Value = 1865
$ws.Cells.Item(116, 9).Value = 1715.5834
$ws.Cells.Item(116, 11).Value = 1715.5834
$ws.Cells.Item(116, 13).Value = 578.4166
$ws.Cells.Item(132, 8).Value = 16981.363
$ws.Cells.Item(132, 9).Value = 1571.9656
$ws.Cells.Item(132, 11).Value = 4715.8968
$ws.Cells.Item(132, 13).Value = -2185.8968
$ws.Cells.Item(136, 8).Value = 3770.5881
$ws.Cells.Item(136, 9).Value = 3846.875
$ws.Cells.Item(136, 10).Value = 2550
$ws.Cells.Item(136, 11).Value = 11540.625
$ws.Cells.Item(136, 12).Value = 7650
$ws.Cells.Item(136, 13).Value = -8990.625
$ws.Cells.Item(136, 14).Value = -12750

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1865
$ws.Cells.Item(3, 9).Value = 1715.5834
$ws.Cells.Item(3, 11).Value = 1715.5834
$ws.Cells.Item(3, 13).Value = -1601.5834
$ws.Cells.Item(134, 8).Value = 3771.697
$ws.Cells.Item(134, 9).Value = 4120.8965
$ws.Cells.Item(134, 10).Value = 1240
$ws.Cells.Item(134, 11).Value = 12362.6895
$ws.Cells.Item(134, 12).Value = 3720
$ws.Cells.Item(134, 13).Value = -9827.6895
$ws.Cells.Item(134, 14).Value = -8790

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3882.8064
$ws.Cells.Item(31, 9).Value = 3209
$ws.Cells.Item(31, 10).Value = 4308.3687
$ws.Cells.Item(31, 11).Value = 3209
$ws.Cells.Item(31, 12).Value = 4308.3687
$ws.Cells.Item(31, 13).Value = -2914
$ws.Cells.Item(31, 14).Value = -4898.3687
$ws.Cells.Item(34, 8).Value = 3882.8064
$ws.Cells.Item(34, 9).Value = 3209
$ws.Cells.Item(34, 10).Value = 4308.3687
$ws.Cells.Item(34, 11).Value = 3209
$ws.Cells.Item(34, 12).Value = 4308.3687
$ws.Cells.Item(34, 13).Value = -3007
$ws.Cells.Item(34, 14).Value = -4712.3687
$ws.Cells.Item(132, 8).Value = 4419.75
$ws.Cells.Item(132, 9).Value = 3633.7778
$ws.Cells.Item(132, 10).Value = 5430.2856
$ws.Cells.Item(132, 11).Value = 10901.3334
$ws.Cells.Item(132, 12).Value = 16290.8568
$ws.Cells.Item(132, 13).Value = -8371.3334
$ws.Cells.Item(132, 14).Value = -21350.8568
$ws.Cells.Item(134, 8).Value = 1064.8611
$ws.Cells.Item(134, 9).Value = 919.43475
$ws.Cells.Item(134, 10).Value = 1322.1538
$ws.Cells.Item(134, 11).Value = 2758.30425
$ws.Cells.Item(134, 12).Value = 3966.4614
$ws.Cells.Item(134, 13).Value = -223.3042500000001
$ws.Cells.Item(134, 14).Value = -9036.4614

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 461.25
$ws.Cells.Item(23, 9).Value = 40
$ws.Cells.Item(23, 10).Value = 714
$ws.Cells.Item(23, 11).Value = 120
$ws.Cells.Item(23, 12).Value = 2142
$ws.Cells.Item(23, 13).Value = 115
$ws.Cells.Item(23, 14).Value = -2612
$ws.Cells.Item(113, 8).Value = 411.66666
$ws.Cells.Item(113, 9).Value = 420
$ws.Cells.Item(113, 10).Value = 410
$ws.Cells.Item(113, 11).Value = 1260
$ws.Cells.Item(113, 12).Value = 1230
$ws.Cells.Item(113, 13).Value = 910
$ws.Cells.Item(113, 14).Value = -5570
$ws.Cells.Item(131, 8).Value = 738.58
$ws.Cells.Item(131, 10).Value = 742.5567
$ws.Cells.Item(131, 12).Value = 2227.6701
$ws.Cells.Item(131, 14).Value = -12307.6701

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 5450.069
$ws.Cells.Item(126, 9).Value = 4441.8823
$ws.Cells.Item(126, 10).Value = 6878.3335
$ws.Cells.Item(126, 11).Value = 13325.6469
$ws.Cells.Item(126, 12).Value = 20635.0005
$ws.Cells.Item(126, 13).Value = -10855.6469
$ws.Cells.Item(126, 14).Value = -25575.0005
$ws.Cells.Item(132, 8).Value = 30666.223
$ws.Cells.Item(132, 9).Value = 3205.4285
$ws.Cells.Item(132, 10).Value = 48141.273
$ws.Cells.Item(132, 11).Value = 9616.2855
$ws.Cells.Item(132, 12).Value = 144423.819
$ws.Cells.Item(132, 13).Value = -7086.2855
$ws.Cells.Item(132, 14).Value = -149483.819

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4417.6665
$ws.Cells.Item(132, 9).Value = 3627
$ws.Cells.Item(132, 10).Value = 5999
$ws.Cells.Item(132, 11).Value = 10881
$ws.Cells.Item(132, 12).Value = 17997
$ws.Cells.Item(132, 13).Value = -8351
$ws.Cells.Item(132, 14).Value = -23057
$ws.Cells.Item(136, 8).Value = 1832.3334
$ws.Cells.Item(136, 9).Value = 1832.3334
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 5497.0002
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -2947.0002
$ws.Cells.Item(136, 14).ClearContents()

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2099.2856
$ws.Cells.Item(132, 9).Value = 1339.2
$ws.Cells.Item(132, 11).Value = 4017.6
$ws.Cells.Item(132, 13).Value = -1487.6
$ws.Cells.Item(136, 8).Value = 31253232
$ws.Cells.Item(136, 9).Value = 47620716
$ws.Cells.Item(136, 10).Value = 6216
$ws.Cells.Item(136, 11).Value = 142862148
$ws.Cells.Item(136, 12).Value = 18648
$ws.Cells.Item(136, 13).Value = -142859598
$ws.Cells.Item(136, 14).Value = -23748
